# Update the cryptocurrency price/volume table with the latest scraped values.
# Values that look like plain numbers (e.g. "0.999") are written with a
# leading apostrophe so Excel stores them as text, matching the original
# inline-string (t="inlineStr") cell type used throughout column D/E instead
# of silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.252.72'
$ws.Range('E2').Value = '  +7.07%  '
$ws.Range('D3').Value = '2.668.16'
$ws.Range('E3').Value = '  +9.25%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''513.62'
$ws.Range('E5').Value = '  +4.67%  '
$ws.Range('D6').Value = '''159.44'
$ws.Range('E6').Value = '  +3.42%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '''0.621'
$ws.Range('E7').Value = '  +1.43%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '''0.998'
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '2.666.62'
$ws.Range('E9').Value = '  +8.93%  '
$ws.Range('E10').Value = '  +9.42%  '
$ws.Range('E11').Value = '  +5.35%  '
$ws.Range('D12').Value = '''0.353'
$ws.Range('E12').Value = '  +4.83%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '3.126.63'
$ws.Range('E14').Value = '  +9.33%  '
$ws.Range('D15').Value = '61.240.98'
$ws.Range('E15').Value = '  +7.07%  '
$ws.Range('D16').Value = '''22.35'
$ws.Range('E16').Value = '  +7.55%  '
$ws.Range('D17').Value = '''0.0000141'
$ws.Range('E17').Value = '  +5.37%  '
$ws.Range('D18').Value = '2.665.31'
$ws.Range('E18').Value = '  +8.84%  '
$ws.Range('D19').Value = '''4.84'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').Value = '''354.30'
$ws.Range('E20').Value = '  +8.09%  '
$ws.Range('E21').Value = '  +6.39%  '
$ws.Range('D22').Value = '''6.19'
$ws.Range('E22').Value = '  +4.30%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '''60.78'
$ws.Range('E24').Value = '  +5.07%  '
$ws.Range('D25').Value = '''0.427'
$ws.Range('E25').Value = '  +4.13%  '
$ws.Range('D26').Value = '2.779.77'
$ws.Range('E26').Value = '  +9.35%  '
$ws.Range('E27').Value = '  +4.62%  '
$ws.Range('D28').Value = '''1.01'
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('D29').Value = '0.0₃0873'
$ws.Range('E29').Value = '  +10.98%  '
$ws.Range('D30').Value = '''7.60'
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '''19.73'
$ws.Range('E32').Value = '  +5.30%  '
$ws.Range('D33').Value = '''157.19'
$ws.Range('E33').Value = '  +5.03%  '
$ws.Range('E34').Value = '  +4.47%  '
$ws.Range('E35').Value = '  +8.90%  '
$ws.Range('E36').Value = '  +11.03%  '
$ws.Range('E37').Value = '  +7.54%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = '''0.887'
$ws.Range('E38').Value = '  +3.67%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.55'
$ws.Range('E39').Value = '  +11.69%  '
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D40').Value = '''0.847'
$ws.Range('E40').Value = '  +32.87%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '''3.81'
$ws.Range('E41').Value = '  +8.38%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''305.51'
$ws.Range('E42').Value = '  +14.38%  '
$ws.Range('D43').Value = '''35.84'
$ws.Range('E43').Value = '  +4.56%  '
$ws.Range('D44').Value = '''0.647'
$ws.Range('E44').Value = '  +7.94%  '
$ws.Range('D45').Value = '''0.0583'
$ws.Range('E45').Value = '  +8.34%  '
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '''20.29'
$ws.Range('E47').Value = '  +15.05%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').Value = '''4.99'
$ws.Range('E49').Value = '  +6.23%  '
$ws.Range('E50').Value = '  +4.28%  '
$ws.Range('D51').Value = '2.029.80'
$ws.Range('E51').Value = '  +8.74%  '
